# Generate Report for handoff
#
# Sets the "Latest Handoff Datetime" (column D) for the most recently
# handed-off file (row 5, 7ec9710a-...) on both the zh-cn and de-de
# language report sheets, recording the new handoff timestamps, and adds
# the corresponding timestamp strings to the shared string table.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-01-15 07:44:33"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-01-15 07:44:42"
